$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A12 to a new testcase id (was a duplicate of A11)
$ws.Range("A12").Value = "ADD_LEAVE_TYPE_TC010"

# Move the active selection to B12
$ws.Range("B12").Select()
